# "Thêm khung Controller BNS" - add a BanNhanSu (BNS) controller account row
# to the Sheet1 account list, and move the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new account row right after the existing data (row 7):
#   UserName = BanNhanSu, Password = 12345, UserType = super
$ws.Range("A7").Value = "BanNhanSu"
$ws.Range("B7").Value = 12345
$ws.Range("C7").Value = "super"

# Update the sheet's active selection/cell as recorded in the saved view.
$ws.Range("G8").Select()
